$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "The applicant has a high score because the projects the applicant has worked on demonstrate expertise in NLP, DL, ML, and computer vision, which are directly relevant to the job description. The applicant has utilized tools such as PyTorch, Python, and Computer Vision in their projects, aligning with the required skills for the job."
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = "The applicant shows a strong willingness to work in Japan and has demonstrated adaptability, communication, teamwork, and a desire to learn new skills and adapt to a new work culture. This aligns well with the soft skills and the company's values."

$ws.Range("C3").Value = "The applicant has a high score due to the extensive experience in NLP, ML, and computer vision using Python, PyTorch, and ReactJS, which aligns with the job requirements."
$ws.Range("E3").Value = "The applicant demonstrates strong adaptability, teamwork, and problem-solving skills, while also expressing a willingness to work in Japan and learn a new language, earning a high score for compatibility with the company's soft skills and willingness to work abroad."
